# Auto-generated cell updates applying the "Updated cryptos list" diff.
# Each entry is Cell, NewValue, IsNumericLooking (needs quote-prefix text protection).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "66.037.53"; Numeric = 0 }
    @{ Cell = "E2"; Value = "  -1.33%  "; Numeric = 0 }
    @{ Cell = "D3"; Value = "3.483.49"; Numeric = 0 }
    @{ Cell = "E3"; Value = "  +0.13%  "; Numeric = 0 }
    @{ Cell = "D4"; Value = "0.999"; Numeric = 1 }
    @{ Cell = "E4"; Value = "  -0.14%  "; Numeric = 0 }
    @{ Cell = "D5"; Value = "603.26"; Numeric = 1 }
    @{ Cell = "E5"; Value = "  +0.36%  "; Numeric = 0 }
    @{ Cell = "D6"; Value = "143.22"; Numeric = 1 }
    @{ Cell = "E6"; Value = "  -3.13%  "; Numeric = 0 }
    @{ Cell = "D7"; Value = "3.481.37"; Numeric = 0 }
    @{ Cell = "E7"; Value = "  +0.12%  "; Numeric = 0 }
    @{ Cell = "D8"; Value = "1.00"; Numeric = 1 }
    @{ Cell = "E8"; Value = "  +0.00%  "; Numeric = 0 }
    @{ Cell = "D9"; Value = "0.475"; Numeric = 1 }
    @{ Cell = "E9"; Value = "  -0.93%  "; Numeric = 0 }
    @{ Cell = "D10"; Value = "8.16"; Numeric = 1 }
    @{ Cell = "E10"; Value = "  +6.94%  "; Numeric = 0 }
    @{ Cell = "D11"; Value = "0.135"; Numeric = 1 }
    @{ Cell = "E11"; Value = "  -4.38%  "; Numeric = 0 }
    @{ Cell = "D12"; Value = "0.413"; Numeric = 1 }
    @{ Cell = "E12"; Value = "  -2.23%  "; Numeric = 0 }
    @{ Cell = "D13"; Value = "4.072.04"; Numeric = 0 }
    @{ Cell = "E13"; Value = "  +0.09%  "; Numeric = 0 }
    @{ Cell = "E14"; Value = "  -3.95%  "; Numeric = 0 }
    @{ Cell = "D15"; Value = "30.36"; Numeric = 1 }
    @{ Cell = "E15"; Value = "  -2.85%  "; Numeric = 0 }
    @{ Cell = "D16"; Value = "3.492.27"; Numeric = 0 }
    @{ Cell = "E16"; Value = "  +0.60%  "; Numeric = 0 }
    @{ Cell = "B17"; Value = "TRON"; Numeric = 0 }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; Numeric = 0 }
    @{ Cell = "D17"; Value = "0.117"; Numeric = 1 }
    @{ Cell = "E17"; Value = "  -0.37%  "; Numeric = 0 }
    @{ Cell = "B18"; Value = "WrappedBTC"; Numeric = 0 }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; Numeric = 0 }
    @{ Cell = "D18"; Value = "66.089.43"; Numeric = 0 }
    @{ Cell = "E18"; Value = "  -1.29%  "; Numeric = 0 }
    @{ Cell = "D19"; Value = "10.44"; Numeric = 1 }
    @{ Cell = "E19"; Value = "  +3.82%  "; Numeric = 0 }
    @{ Cell = "D20"; Value = "6.19"; Numeric = 1 }
    @{ Cell = "E20"; Value = "  -3.29%  "; Numeric = 0 }
    @{ Cell = "D21"; Value = "14.73"; Numeric = 1 }
    @{ Cell = "E21"; Value = "  -3.52%  "; Numeric = 0 }
    @{ Cell = "D22"; Value = "420.63"; Numeric = 1 }
    @{ Cell = "E22"; Value = "  -2.98%  "; Numeric = 0 }
    @{ Cell = "D23"; Value = "0.588"; Numeric = 1 }
    @{ Cell = "E23"; Value = "  -2.79%  "; Numeric = 0 }
    @{ Cell = "D24"; Value = "77.63"; Numeric = 1 }
    @{ Cell = "E24"; Value = "  -1.71%  "; Numeric = 0 }
    @{ Cell = "D26"; Value = "0.0000116"; Numeric = 1 }
    @{ Cell = "E26"; Value = "  -3.13%  "; Numeric = 0 }
    @{ Cell = "D27"; Value = "9.45"; Numeric = 1 }
    @{ Cell = "E27"; Value = "  -3.45%  "; Numeric = 0 }
    @{ Cell = "D28"; Value = "7.97"; Numeric = 1 }
    @{ Cell = "E28"; Value = "  -4.83%  "; Numeric = 0 }
    @{ Cell = "D29"; Value = "2.46"; Numeric = 1 }
    @{ Cell = "E29"; Value = "  -0.70%  "; Numeric = 0 }
    @{ Cell = "E30"; Value = "  -0.45%  "; Numeric = 0 }
    @{ Cell = "D31"; Value = "0.162"; Numeric = 1 }
    @{ Cell = "E31"; Value = "  -3.34%  "; Numeric = 0 }
    @{ Cell = "D32"; Value = "1.48"; Numeric = 1 }
    @{ Cell = "E32"; Value = "  -6.02%  "; Numeric = 0 }
    @{ Cell = "D33"; Value = "25.14"; Numeric = 1 }
    @{ Cell = "E33"; Value = "  -0.54%  "; Numeric = 0 }
    @{ Cell = "D34"; Value = "3.477.53"; Numeric = 0 }
    @{ Cell = "E34"; Value = "  +0.19%  "; Numeric = 0 }
    @{ Cell = "E35"; Value = "  -0.08%  "; Numeric = 0 }
    @{ Cell = "D36"; Value = "1.71"; Numeric = 1 }
    @{ Cell = "E36"; Value = "  -4.63%  "; Numeric = 0 }
    @{ Cell = "D37"; Value = "5.57"; Numeric = 1 }
    @{ Cell = "E37"; Value = "  -5.75%  "; Numeric = 0 }
    @{ Cell = "D38"; Value = "7.65"; Numeric = 1 }
    @{ Cell = "E38"; Value = "  -2.88%  "; Numeric = 0 }
    @{ Cell = "D39"; Value = "0.999"; Numeric = 1 }
    @{ Cell = "E39"; Value = "  -0.04%  "; Numeric = 0 }
    @{ Cell = "D40"; Value = "170.39"; Numeric = 1 }
    @{ Cell = "E40"; Value = "  -1.83%  "; Numeric = 0 }
    @{ Cell = "D41"; Value = "0.0869"; Numeric = 1 }
    @{ Cell = "E41"; Value = "  -1.54%  "; Numeric = 0 }
    @{ Cell = "D42"; Value = "0.892"; Numeric = 1 }
    @{ Cell = "E42"; Value = "  -0.18%  "; Numeric = 0 }
    @{ Cell = "D43"; Value = "5.13"; Numeric = 1 }
    @{ Cell = "E43"; Value = "  -4.66%  "; Numeric = 0 }
    @{ Cell = "D44"; Value = "1.91"; Numeric = 1 }
    @{ Cell = "E44"; Value = "  -8.17%  "; Numeric = 0 }
    @{ Cell = "D45"; Value = "45.68"; Numeric = 1 }
    @{ Cell = "E45"; Value = "  -1.60%  "; Numeric = 0 }
    @{ Cell = "D46"; Value = "26.18"; Numeric = 1 }
    @{ Cell = "E46"; Value = "  -9.38%  "; Numeric = 0 }
    @{ Cell = "E47"; Value = "  -2.64%  "; Numeric = 0 }
    @{ Cell = "B48"; Value = "Cosmos"; Numeric = 0 }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; Numeric = 0 }
    @{ Cell = "D48"; Value = "7.13"; Numeric = 1 }
    @{ Cell = "E48"; Value = "  -4.14%  "; Numeric = 0 }
    @{ Cell = "B49"; Value = "dogwifhat"; Numeric = 0 }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; Numeric = 0 }
    @{ Cell = "D49"; Value = "2.35"; Numeric = 1 }
    @{ Cell = "E49"; Value = "  -2.15%  "; Numeric = 0 }
    @{ Cell = "D50"; Value = "0.936"; Numeric = 1 }
    @{ Cell = "E50"; Value = "  -3.74%  "; Numeric = 0 }
    @{ Cell = "D51"; Value = "0.236"; Numeric = 1 }
    @{ Cell = "E51"; Value = "  -3.17%  "; Numeric = 0 }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Numeric -eq 1) {
        # Prefix with an apostrophe so Excel stores it as text (preserving trailing
        # zeros / exact formatting) instead of converting it to a number, then reset
        # the cell style back to Normal so no stray style index remains attached.
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
